$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains its text semantics (values like "1.007" would otherwise
# be auto-converted to numbers by Excel). Apply Text number format to the whole
# price column first so individual assignments below keep their exact string form.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value2 = "27.067.18"
$ws.Range("E2").Value2 = "  -1.75%  "

$ws.Range("D3").Value2 = "1.794.90"
$ws.Range("E3").Value2 = "  -2.17%  "

$ws.Range("D4").Value2 = "1.007"
$ws.Range("E4").Value2 = "  +0.18%  "

$ws.Range("D5").Value2 = "1.007"
$ws.Range("E5").Value2 = "  +0.31%  "

$ws.Range("D6").Value2 = "308.63"
$ws.Range("E6").Value2 = "  -1.49%  "

$ws.Range("D7").Value2 = "0.4166"
$ws.Range("E7").Value2 = "  -1.75%  "

$ws.Range("D8").Value2 = "0.3545"
$ws.Range("E8").Value2 = "  -3.16%  "

$ws.Range("D9").Value2 = "0.07043"
$ws.Range("E9").Value2 = "  -2.49%  "

$ws.Range("D10").Value2 = "0.8442"
$ws.Range("E10").Value2 = "  -2.41%  "

$ws.Range("D11").Value2 = "20.12"
$ws.Range("E11").Value2 = "  -3.08%  "

$ws.Range("D12").Value2 = "1.869.37"
$ws.Range("E12").Value2 = "  +0.43%  "

$ws.Range("D13").Value2 = "5.251"
$ws.Range("E13").Value2 = "  -2.30%  "

$ws.Range("D14").Value2 = "6.330"
$ws.Range("E14").Value2 = "  -2.62%  "

$ws.Range("D15").Value2 = "0.06840"
$ws.Range("E15").Value2 = "  -1.68%  "

$ws.Range("D16").Value2 = "1.010"
$ws.Range("E16").Value2 = "  +0.32%  "

$ws.Range("D17").Value2 = "79.65"
$ws.Range("E17").Value2 = "  -0.15%  "

$ws.Range("D18").Value2 = "0.000008716"
$ws.Range("E18").Value2 = "  -2.89%  "

$ws.Range("D19").Value2 = "1.008"
$ws.Range("E19").Value2 = "  +0.49%  "

$ws.Range("D20").Value2 = "15.04"
$ws.Range("E20").Value2 = "  -2.19%  "

$ws.Range("D21").Value2 = "27.590.51"
$ws.Range("E21").Value2 = "  -0.32%  "

$ws.Range("D22").Value2 = "5.038"
$ws.Range("E22").Value2 = "  +0.05%  "

$ws.Range("D23").Value2 = "10.72"
$ws.Range("E23").Value2 = "  -0.67%  "

$ws.Range("D24").Value2 = "2.114.98"
$ws.Range("E24").Value2 = "  +0.56%  "

$ws.Range("D25").Value2 = "1.957"
$ws.Range("E25").Value2 = "  -0.39%  "

$ws.Range("D26").Value2 = "152.75"
$ws.Range("E26").Value2 = "  -0.91%  "

$ws.Range("D27").Value2 = "18.14"
$ws.Range("E27").Value2 = "  -1.32%  "

$ws.Range("D28").Value2 = "5.015"
$ws.Range("E28").Value2 = "  -4.27%  "

$ws.Range("D29").Value2 = "112.22"
$ws.Range("E29").Value2 = "  -2.09%  "

$ws.Range("D30").Value2 = "1.646"
$ws.Range("E30").Value2 = "  -9.40%  "

$ws.Range("D31").Value2 = "0.08857"
$ws.Range("E31").Value2 = "  -0.11%  "

$ws.Range("D32").Value2 = "0.7238"
$ws.Range("E32").Value2 = "  -6.16%  "

$ws.Range("D33").Value2 = "2.876"
$ws.Range("E33").Value2 = "  -2.43%  "

$ws.Range("D34").Value2 = "4.345"
$ws.Range("E34").Value2 = "  -4.49%  "

$ws.Range("D35").Value2 = "1.007"
$ws.Range("E35").Value2 = "  +0.37%  "

$ws.Range("D36").Value2 = "1.080"
$ws.Range("E36").Value2 = "  -5.77%  "

$ws.Range("D37").Value2 = "1.074"
$ws.Range("E37").Value2 = "  -2.05%  "

$ws.Range("D38").Value2 = "0.05101"

$ws.Range("D39").Value2 = "0.01890"
$ws.Range("E39").Value2 = "  -2.68%  "

$ws.Range("D40").Value2 = "0.4933"
$ws.Range("E40").Value2 = "  -3.52%  "

$ws.Range("D41").Value2 = "0.1612"
$ws.Range("E41").Value2 = "  -2.11%  "

$ws.Range("D42").Value2 = "2.631"
$ws.Range("E42").Value2 = "  -6.96%  "

$ws.Range("D43").Value2 = "6.177"
$ws.Range("E43").Value2 = "  -9.29%  "

$ws.Range("D44").Value2 = "8.034"
$ws.Range("E44").Value2 = "  -4.94%  "

$ws.Range("D45").Value2 = "1.007"
$ws.Range("E45").Value2 = "  +0.45%  "

$ws.Range("D46").Value2 = "10.23"
$ws.Range("E46").Value2 = "  -2.16%  "

$ws.Range("D47").Value2 = "104.03"
$ws.Range("E47").Value2 = "  -1.59%  "

$ws.Range("D48").Value2 = "0.06314"
$ws.Range("E48").Value2 = "  -3.37%  "

$ws.Range("D49").Value2 = "0.4528"
$ws.Range("E49").Value2 = "  -3.45%  "

$ws.Range("D50").Value2 = "1.583"
$ws.Range("E50").Value2 = "  -2.35%  "

$ws.Range("D51").Value2 = "62.09"
$ws.Range("E51").Value2 = "  -3.18%  "
